$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate row 2 with the GameServer_1 entry.
# Columns: A=ID, B=ServerID, C=Name, D=MaxOnline, E=CpuCount, F=IP, G=Port
$ws.Range("A2").Value = "GameServer_1"
$ws.Range("B2").Value = "000104001"

# C2 and F2 need the same "text" number format (style) the existing
# A2/B2 cells already carry, since they're brand-new cells in this row.
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "GameServer_1"

$ws.Range("D2").Value = 5000
$ws.Range("E2").Value = 1

$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "127.0.0.1"

$ws.Range("G2").Value = 4001

# Move the data validation list down so it starts at F3 (below the new
# data row) instead of F2.
$ws.Range("F2:F1048576").Validation.Delete()
$ws.Range("F3:F1048576").Validation.Add(3, 1, 0, '"TRUE,FALSE"')

# Update the active selection to G3.
$ws.Range("G3").Select()
